# "final version of STOSS-Matlab"
# Update the three input values on Sheet1 and move the active selection
# to match where the user was last working (B7:C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Initial investment amount: 500,000 -> 1,000,000
$ws.Range("B4").Value = 1000000

# Third parameter: 1 -> 0
$ws.Range("B6").Value = 0

# Fourth parameter: 26 -> 27
$ws.Range("B7").Value = 27

# Move/extend the selection to the merged B7:C7 cell (was B4:C4)
$ws.Range("B7:C7").Select()
